$wb = $excel.ActiveWorkbook

# --- Rename header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
# Duplicate the "Weekly Quantity" sheet so that sheet-level formatting
# (sheetPr, pageMargins, header/date cell styles) matches the rest of the
# workbook, then clear out its data and fill in the forecast contents.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWeekly.Copy([System.Reflection.Missing]::Value, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "PO Forecast"

# Remove the copied data rows, keep only the header row.
$newSheet.Rows("2:5").ClearContents()

# Extend the header style (bold, border, center/top alignment) to C1:D1.
$newSheet.Range("A1:B1").Copy($newSheet.Range("C1:D1"))

# Extend the date number-format style down column A through row 13.
$newSheet.Range("A2").Copy($newSheet.Range("A2:A13"))

# --- Header values ---
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$newSheet.Range("A2").Value = 44934.99999999999
$newSheet.Range("B2").Value = 7
$newSheet.Range("C2").Value = 4.155927414079255
$newSheet.Range("D2").Value = 10.34492014781487

$newSheet.Range("A3").Value = 44941.99999999999
$newSheet.Range("B3").Value = 8
$newSheet.Range("C3").Value = 5.059323584771657
$newSheet.Range("D3").Value = 10.86887277883929

$newSheet.Range("A4").Value = 44955.99999999999
$newSheet.Range("B4").Value = 9
$newSheet.Range("C4").Value = 6.120740409621619
$newSheet.Range("D4").Value = 12.24242447460855

$newSheet.Range("A5").Value = 44962.99999999999
$newSheet.Range("B5").Value = 10
$newSheet.Range("C5").Value = 6.360056504078687
$newSheet.Range("D5").Value = 12.57996630174155

$newSheet.Range("A6").Value = 44969.99999999999
$newSheet.Range("B6").Value = 10
$newSheet.Range("C6").Value = 7.180337430604695
$newSheet.Range("D6").Value = 13.12791546209721

$newSheet.Range("A7").Value = 44976.99999999999
$newSheet.Range("B7").Value = 11
$newSheet.Range("C7").Value = 7.843374676210944
$newSheet.Range("D7").Value = 13.71265499214711

$newSheet.Range("A8").Value = 44983.99999999999
$newSheet.Range("B8").Value = 11
$newSheet.Range("C8").Value = 8.35572680725422
$newSheet.Range("D8").Value = 14.75820971205549

$newSheet.Range("A9").Value = 44990.99999999999
$newSheet.Range("B9").Value = 12
$newSheet.Range("C9").Value = 9.339618819863732
$newSheet.Range("D9").Value = 14.96528806839969

$newSheet.Range("A10").Value = 44997.99999999999
$newSheet.Range("B10").Value = 13
$newSheet.Range("C10").Value = 9.46235664308516
$newSheet.Range("D10").Value = 15.75275219951453

$newSheet.Range("A11").Value = 45004.99999999999
$newSheet.Range("B11").Value = 13
$newSheet.Range("C11").Value = 10.1109050186153
$newSheet.Range("D11").Value = 16.33935694482792

$newSheet.Range("A12").Value = 45011.99999999999
$newSheet.Range("B12").Value = 14
$newSheet.Range("C12").Value = 10.78210850761525
$newSheet.Range("D12").Value = 16.99760150013527

$newSheet.Range("A13").Value = 45018.99999999999
$newSheet.Range("B13").Value = 14
$newSheet.Range("C13").Value = 11.67058044115921
$newSheet.Range("D13").Value = 17.32131574978375
